$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Rps19"
$ws.Cells.Item(2, 3).Value = "C5ar1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 82.22439799999999
$ws.Cells.Item(2, 8).Value = 246.673194
$ws.Cells.Item(2, 9).Value = 0.16435464576988
$ws.Cells.Item(2, 10).Value = 0.16435464576988
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.3893403333333333
$ws.Cells.Item(2, 14).Value = 1.168021
$ws.Cells.Item(2, 15).Value = 0.01344665141573655
$ws.Cells.Item(2, 16).Value = 0.01344665141573655
$ws.Cells.Item(2, 17).Value = 32.01327452545267
$ws.Cells.Item(2, 18).Value = 288.119470729074
$ws.Cells.Item(2, 19).Value = 0.002210019630224437
$ws.Cells.Item(2, 20).Value = 0.002210019630224437

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Rps19"
$ws.Cells.Item(3, 3).Value = "C5ar1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 82.22439799999999
$ws.Cells.Item(3, 8).Value = 246.673194
$ws.Cells.Item(3, 9).Value = 0.16435464576988
$ws.Cells.Item(3, 10).Value = 0.16435464576988
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.193104333333333
$ws.Cells.Item(3, 14).Value = 3.579313
$ws.Cells.Item(3, 15).Value = 0.04120625760907917
$ws.Cells.Item(3, 16).Value = 0.04120625760907917
$ws.Cells.Item(3, 17).Value = 98.10228555952466
$ws.Cells.Item(3, 18).Value = 882.9205700357219
$ws.Cells.Item(3, 19).Value = 0.006772439872842629
$ws.Cells.Item(3, 20).Value = 0.006772439872842628

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Rps19"
$ws.Cells.Item(4, 3).Value = "C5ar1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 82.22439799999999
$ws.Cells.Item(4, 8).Value = 246.673194
$ws.Cells.Item(4, 9).Value = 0.16435464576988
$ws.Cells.Item(4, 10).Value = 0.16435464576988
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.087271
$ws.Cells.Item(4, 14).Value = 0.261813
$ws.Cells.Item(4, 15).Value = 0.003014079496094877
$ws.Cells.Item(4, 16).Value = 0.003014079496094878
$ws.Cells.Item(4, 17).Value = 7.175805437857999
$ws.Cells.Item(4, 18).Value = 64.582248940722
$ws.Cells.Item(4, 19).Value = 0.000495377967902932
$ws.Cells.Item(4, 20).Value = 0.000495377967902932

# Row 5: ECs -> Resolving-Mac
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Rps19"
$ws.Cells.Item(5, 3).Value = "C5ar1"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 82.22439799999999
$ws.Cells.Item(5, 8).Value = 246.673194
$ws.Cells.Item(5, 9).Value = 0.16435464576988
$ws.Cells.Item(5, 10).Value = 0.16435464576988
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 27.28472966666666
$ws.Cells.Item(5, 14).Value = 81.85418899999999
$ws.Cells.Item(5, 15).Value = 0.9423330114790893
$ws.Cells.Item(5, 16).Value = 0.9423330114790894
$ws.Cells.Item(5, 17).Value = 2243.470471434407
$ws.Cells.Item(5, 18).Value = 20191.23424290966
$ws.Cells.Item(5, 19).Value = 0.15487680829891
$ws.Cells.Item(5, 20).Value = 0.15487680829891

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Rps19"
$ws.Cells.Item(6, 3).Value = "C5ar1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 123.930687
$ws.Cells.Item(6, 8).Value = 371.792061
$ws.Cells.Item(6, 9).Value = 0.2477194683979672
$ws.Cells.Item(6, 10).Value = 0.2477194683979671
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.3893403333333333
$ws.Cells.Item(6, 14).Value = 1.168021
$ws.Cells.Item(6, 15).Value = 0.01344665141573655
$ws.Cells.Item(6, 16).Value = 0.01344665141573655
$ws.Cells.Item(6, 17).Value = 48.251214986809
$ws.Cells.Item(6, 18).Value = 434.260934881281
$ws.Cells.Item(6, 19).Value = 0.003330997340439032
$ws.Cells.Item(6, 20).Value = 0.003330997340439032

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Rps19"
$ws.Cells.Item(7, 3).Value = "C5ar1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 123.930687
$ws.Cells.Item(7, 8).Value = 371.792061
$ws.Cells.Item(7, 9).Value = 0.2477194683979672
$ws.Cells.Item(7, 10).Value = 0.2477194683979671
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.193104333333333
$ws.Cells.Item(7, 14).Value = 3.579313
$ws.Cells.Item(7, 15).Value = 0.04120625760907917
$ws.Cells.Item(7, 16).Value = 0.04120625760907917
$ws.Cells.Item(7, 17).Value = 147.862239692677
$ws.Cells.Item(7, 18).Value = 1330.760157234093
$ws.Cells.Item(7, 19).Value = 0.01020759222959078
$ws.Cells.Item(7, 20).Value = 0.01020759222959078

# Row 8: FAPs -> MuSCs
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Rps19"
$ws.Cells.Item(8, 3).Value = "C5ar1"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 123.930687
$ws.Cells.Item(8, 8).Value = 371.792061
$ws.Cells.Item(8, 9).Value = 0.2477194683979672
$ws.Cells.Item(8, 10).Value = 0.2477194683979671
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.087271
$ws.Cells.Item(8, 14).Value = 0.261813
$ws.Cells.Item(8, 15).Value = 0.003014079496094877
$ws.Cells.Item(8, 16).Value = 0.003014079496094878
$ws.Cells.Item(8, 17).Value = 10.815554985177
$ws.Cells.Item(8, 18).Value = 97.339994866593
$ws.Cells.Item(8, 19).Value = 0.0007466461704818357
$ws.Cells.Item(8, 20).Value = 0.0007466461704818357

# Row 9: FAPs -> Resolving-Mac
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Rps19"
$ws.Cells.Item(9, 3).Value = "C5ar1"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 123.930687
$ws.Cells.Item(9, 8).Value = 371.792061
$ws.Cells.Item(9, 9).Value = 0.2477194683979672
$ws.Cells.Item(9, 10).Value = 0.2477194683979671
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 27.28472966666666
$ws.Cells.Item(9, 14).Value = 81.85418899999999
$ws.Cells.Item(9, 15).Value = 0.9423330114790893
$ws.Cells.Item(9, 16).Value = 0.9423330114790894
$ws.Cells.Item(9, 17).Value = 3381.41529219928
$ws.Cells.Item(9, 18).Value = 30432.73762979353
$ws.Cells.Item(9, 19).Value = 0.2334342326574555
$ws.Cells.Item(9, 20).Value = 0.2334342326574555

# Row 10: MuSCs -> ECs
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Rps19"
$ws.Cells.Item(10, 3).Value = "C5ar1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 183.1085766666667
$ws.Cells.Item(10, 8).Value = 549.32573
$ws.Cells.Item(10, 9).Value = 0.3660074866766057
$ws.Cells.Item(10, 10).Value = 0.3660074866766056
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.3893403333333333
$ws.Cells.Item(10, 14).Value = 1.168021
$ws.Cells.Item(10, 15).Value = 0.01344665141573655
$ws.Cells.Item(10, 16).Value = 0.01344665141573655
$ws.Cells.Item(10, 17).Value = 71.29155427559222
$ws.Cells.Item(10, 18).Value = 641.62398848033
$ws.Cells.Item(10, 19).Value = 0.004921575088890157
$ws.Cells.Item(10, 20).Value = 0.004921575088890157

# Row 11: MuSCs -> FAPs
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Rps19"
$ws.Cells.Item(11, 3).Value = "C5ar1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 183.1085766666667
$ws.Cells.Item(11, 8).Value = 549.32573
$ws.Cells.Item(11, 9).Value = 0.3660074866766057
$ws.Cells.Item(11, 10).Value = 0.3660074866766056
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.193104333333333
$ws.Cells.Item(11, 14).Value = 3.579313
$ws.Cells.Item(11, 15).Value = 0.04120625760907917
$ws.Cells.Item(11, 16).Value = 0.04120625760907917
$ws.Cells.Item(11, 17).Value = 218.4676362914989
$ws.Cells.Item(11, 18).Value = 1966.20872662349
$ws.Cells.Item(11, 19).Value = 0.01508179878284782
$ws.Cells.Item(11, 20).Value = 0.01508179878284782

# Row 12: MuSCs -> MuSCs
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Rps19"
$ws.Cells.Item(12, 3).Value = "C5ar1"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 183.1085766666667
$ws.Cells.Item(12, 8).Value = 549.32573
$ws.Cells.Item(12, 9).Value = 0.3660074866766057
$ws.Cells.Item(12, 10).Value = 0.3660074866766056
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.087271
$ws.Cells.Item(12, 14).Value = 0.261813
$ws.Cells.Item(12, 15).Value = 0.003014079496094877
$ws.Cells.Item(12, 16).Value = 0.003014079496094878
$ws.Cells.Item(12, 17).Value = 15.98006859427667
$ws.Cells.Item(12, 18).Value = 143.82061734849
$ws.Cells.Item(12, 19).Value = 0.001103175661009176
$ws.Cells.Item(12, 20).Value = 0.001103175661009176

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Rps19"
$ws.Cells.Item(13, 3).Value = "C5ar1"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 183.1085766666667
$ws.Cells.Item(13, 8).Value = 549.32573
$ws.Cells.Item(13, 9).Value = 0.3660074866766057
$ws.Cells.Item(13, 10).Value = 0.3660074866766056
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 27.28472966666666
$ws.Cells.Item(13, 14).Value = 81.85418899999999
$ws.Cells.Item(13, 15).Value = 0.9423330114790893
$ws.Cells.Item(13, 16).Value = 0.9423330114790894
$ws.Cells.Item(13, 17).Value = 4996.068013998107
$ws.Cells.Item(13, 18).Value = 44964.61212598297
$ws.Cells.Item(13, 19).Value = 0.3449009371438584
$ws.Cells.Item(13, 20).Value = 0.3449009371438584

# Row 14: Resolving-Mac -> ECs
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Rps19"
$ws.Cells.Item(14, 3).Value = "C5ar1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 111.0227623333333
$ws.Cells.Item(14, 8).Value = 333.0682870000001
$ws.Cells.Item(14, 9).Value = 0.2219183991555473
$ws.Cells.Item(14, 10).Value = 0.2219183991555472
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.3893403333333333
$ws.Cells.Item(14, 14).Value = 1.168021
$ws.Cells.Item(14, 15).Value = 0.01344665141573655
$ws.Cells.Item(14, 16).Value = 0.01344665141573655
$ws.Cells.Item(14, 17).Value = 43.22563929444745
$ws.Cells.Item(14, 18).Value = 389.0307536500271
$ws.Cells.Item(14, 19).Value = 0.002984059356182929
$ws.Cells.Item(14, 20).Value = 0.002984059356182929

# Row 15: Resolving-Mac -> FAPs
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Rps19"
$ws.Cells.Item(15, 3).Value = "C5ar1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 111.0227623333333
$ws.Cells.Item(15, 8).Value = 333.0682870000001
$ws.Cells.Item(15, 9).Value = 0.2219183991555473
$ws.Cells.Item(15, 10).Value = 0.2219183991555472
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 1.193104333333333
$ws.Cells.Item(15, 14).Value = 3.579313
$ws.Cells.Item(15, 15).Value = 0.04120625760907917
$ws.Cells.Item(15, 16).Value = 0.04120625760907917
$ws.Cells.Item(15, 17).Value = 132.4617388385368
$ws.Cells.Item(15, 18).Value = 1192.155649546831
$ws.Cells.Item(15, 19).Value = 0.009144426723797937
$ws.Cells.Item(15, 20).Value = 0.009144426723797935

# Row 16: Resolving-Mac -> MuSCs
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Rps19"
$ws.Cells.Item(16, 3).Value = "C5ar1"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 111.0227623333333
$ws.Cells.Item(16, 8).Value = 333.0682870000001
$ws.Cells.Item(16, 9).Value = 0.2219183991555473
$ws.Cells.Item(16, 10).Value = 0.2219183991555472
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.087271
$ws.Cells.Item(16, 14).Value = 0.261813
$ws.Cells.Item(16, 15).Value = 0.003014079496094877
$ws.Cells.Item(16, 16).Value = 0.003014079496094878
$ws.Cells.Item(16, 17).Value = 9.689067491592334
$ws.Cells.Item(16, 18).Value = 87.20160742433102
$ws.Cells.Item(16, 19).Value = 0.0006688796967009338
$ws.Cells.Item(16, 20).Value = 0.0006688796967009338

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Rps19"
$ws.Cells.Item(17, 3).Value = "C5ar1"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 111.0227623333333
$ws.Cells.Item(17, 8).Value = 333.0682870000001
$ws.Cells.Item(17, 9).Value = 0.2219183991555473
$ws.Cells.Item(17, 10).Value = 0.2219183991555472
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 27.28472966666666
$ws.Cells.Item(17, 14).Value = 81.85418899999999
$ws.Cells.Item(17, 15).Value = 0.9423330114790893
$ws.Cells.Item(17, 16).Value = 0.9423330114790894
$ws.Cells.Item(17, 17).Value = 3029.226057111583
$ws.Cells.Item(17, 18).Value = 27263.03451400424
$ws.Cells.Item(17, 19).Value = 0.2091210333788655
$ws.Cells.Item(17, 20).Value = 0.2091210333788655
